$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers in row 1: "<name>_old" -> "<name>_FV2404" and
#    "<name>_new" -> "<name>_FV2410" (column K stays "diff").
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the header row + data into an Excel Table ("Table1") with an
#    autofilter, spanning the whole used range A1:U68.
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$tableRange = $ws.Range("A1:U$lastRow")
$tbl = $ws.ListObjects.Add(1, $tableRange, [Type]::Missing, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
